$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.786.54'
$ws.Range("E2").Value = '  +1.70%  '

$ws.Range("D3").Value = '2.799.55'
$ws.Range("E3").Value = '  +1.69%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '350.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.06%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '112.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.49%  '

$ws.Range("E7").Value = '  +2.39%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("E9").Value = '  +6.99%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.23'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.61%  '

$ws.Range("E11").Value = '  -0.50%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0838'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.93%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.78'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.18%  '

$ws.Range("D15").Value = '3.237.54'
$ws.Range("E15").Value = '  +1.72%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.967'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.56%  '

$ws.Range("D17").Value = '2.812.08'
$ws.Range("E17").Value = '  +3.07%  '

$ws.Range("D18").Value = '51.769.10'
$ws.Range("E18").Value = '  +1.69%  '

$ws.Range("E19").Value = '  +10.52%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.63'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.91%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.55'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.40%  '

$ws.Range("D22").Value = '0.0₃0974'
$ws.Range("E22").Value = '  +1.89%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.41'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.48%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '269.42'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.99%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.75'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.35%  '

$ws.Range("E26").Value = '  -0.10%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.20%  '

$ws.Range("E28").Value = '  +0.51%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '38.86'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +14.05%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.45'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.75%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.28'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.44%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '52.79'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.29%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.12'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.93%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0909'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +9.97%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0455'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.95%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.65'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.66%  '

$ws.Range("E37").Value = '  -0.09%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.88'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.47%  '

$ws.Range("E39").Value = '  +1.56%  '

$ws.Range("E41").Value = '  +2.36%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.52'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.57%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '122.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.27%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.94'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.90%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.53'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.58%  '

$ws.Range("E47").Value = '  +8.97%  '

$ws.Range("D48").Value = '2.121.51'
$ws.Range("E48").Value = '  +1.97%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.984'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.65%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.49'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.05%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.222'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +17.34%  '
